# "finish dev of treatSoldiers"
# Adds 5 new columns (T:X) to the "normal" sheet holding the treat-soldiers
# (healing) costs: INT_treatWood, INT_treatStone, INT_treatIron,
# INT_treatFood, INT_treatTime - mirroring the existing recruit cost
# columns O:S (INT_wood, INT_stone, INT_iron, INT_food, INT_recruitTime).

$wb = $excel.ActiveWorkbook
$wsNormal  = $wb.Worksheets.Item("normal")
$wsSpecial = $wb.Worksheets.Item("special")

# New column headers (row 1), same phrasing style as the existing ones.
$headers = @("INT_treatWood", "INT_treatStone", "INT_treatIron", "INT_treatFood", "INT_treatTime")
$newCols = @(20, 21, 22, 23, 24)   # T, U, V, W, X

# Copy the formatting (number format / borders / alignment) from the
# existing O1:S25 block onto the new T1:X25 block before filling in values,
# so the new cells get the same visual styling as the analogous
# wood/stone/iron/food/time columns.
$wsNormal.Range("O1:S25").Copy() | Out-Null
$wsNormal.Range("T1:X25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new header row.
for ($i = 0; $i -lt $headers.Count; $i++) {
    $wsNormal.Cells.Item(1, $newCols[$i]).Value = $headers[$i]
}

# Row data: rowNumber, treatWood, treatStone, treatIron, treatFood, treatTime
$rows = @(
    @(2,  0,   167, 100, 67,  10),
    @(3,  0,   284, 170, 114, 10),
    @(4,  0,   400, 240, 160, 10),
    @(5,  0,   147, 147, 74,  12),
    @(6,  0,   250, 250, 125, 12),
    @(7,  0,   352, 352, 176, 12),
    @(8,  175, 105, 0,   70,  17),
    @(9,  298, 179, 0,   119, 17),
    @(10, 420, 252, 0,   168, 17),
    @(11, 175, 105, 0,   70,  19),
    @(12, 298, 179, 0,   119, 19),
    @(13, 420, 252, 0,   168, 19),
    @(14, 200, 0,   334, 134, 37),
    @(15, 340, 0,   567, 227, 37),
    @(16, 480, 0,   800, 320, 37),
    @(17, 294, 0,   294, 147, 40),
    @(18, 499, 0,   499, 250, 40),
    @(19, 704, 0,   704, 352, 40),
    @(20, 400, 400, 134, 267, 60),
    @(21, 680, 680, 227, 454, 60),
    @(22, 960, 960, 320, 640, 60),
    @(23, 367, 367, 367, 367, 64),
    @(24, 624, 624, 624, 624, 64),
    @(25, 880, 880, 880, 880, 64)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    for ($i = 0; $i -lt $newCols.Count; $i++) {
        $wsNormal.Cells.Item($r, $newCols[$i]).Value = $entry[$i + 1]
    }
}

# The author finished editing on the "normal" sheet, so it becomes the
# active tab/sheet (instead of "special"), and leaves the selection on
# AC31 there; the "special" sheet selection stays as it was (Q6).
$wsNormal.Activate() | Out-Null
$wsNormal.Range("AC31").Select() | Out-Null
